$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.487.17'
$ws.Cells.Item(2, 5).Value = '  +0.13%  '
$ws.Cells.Item(3, 4).Value = '3.816.71'
$ws.Cells.Item(3, 5).Value = '  +3.49%  '
$ws.Cells.Item(4, 5).Value = '  +0.35%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '410.10'
$ws.Cells.Item(5, 5).Value = '  -2.02%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '130.06'
$ws.Cells.Item(6, 5).Value = '  -0.01%  '
$ws.Cells.Item(7, 4).Value = '3.802.69'
$ws.Cells.Item(7, 5).Value = '  +3.32%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.612'
$ws.Cells.Item(8, 5).Value = '  -4.59%  '
$ws.Cells.Item(9, 5).Value = '  +0.09%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.729'
$ws.Cells.Item(10, 5).Value = '  -6.30%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.168'
$ws.Cells.Item(11, 5).Value = '  -5.59%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000361'
$ws.Cells.Item(12, 5).Value = '  -8.27%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '40.85'
$ws.Cells.Item(13, 5).Value = '  -5.22%  '
$ws.Cells.Item(14, 4).Value = '4.434.83'
$ws.Cells.Item(14, 5).Value = '  +4.03%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '9.95'
$ws.Cells.Item(15, 5).Value = '  -6.18%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '15.58'
$ws.Cells.Item(16, 5).Value = '  +17.18%  '
$ws.Cells.Item(17, 5).Value = '  -1.14%  '
$ws.Cells.Item(18, 4).Value = '3.815.67'
$ws.Cells.Item(18, 5).Value = '  +3.51%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '19.52'
$ws.Cells.Item(19, 5).Value = '  -5.17%  '
$ws.Cells.Item(20, 4).Value = '67.099.32'
$ws.Cells.Item(20, 5).Value = '  +1.04%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.06'
$ws.Cells.Item(21, 5).Value = '  -5.51%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '408.99'
$ws.Cells.Item(22, 5).Value = '  -7.95%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '14.37'
$ws.Cells.Item(23, 5).Value = '  -13.27%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '84.88'
$ws.Cells.Item(24, 5).Value = '  -5.58%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '3.04'
$ws.Cells.Item(25, 5).Value = '  -3.23%  '
$ws.Cells.Item(26, 2).Value = 'EthereumClassic'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '36.25'
$ws.Cells.Item(26, 5).Value = '  -2.48%  '
$ws.Cells.Item(27, 2).Value = 'LEO'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '5.65'
$ws.Cells.Item(27, 5).Value = '  +12.27%  '
$ws.Cells.Item(28, 5).Value = '  -6.12%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.34'
$ws.Cells.Item(29, 5).Value = '  -8.30%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '684.73'
$ws.Cells.Item(30, 5).Value = '  +5.53%  '
$ws.Cells.Item(31, 2).Value = 'Cosmos'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '12.39'
$ws.Cells.Item(31, 5).Value = '  -2.54%  '
$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.121'
$ws.Cells.Item(32, 5).Value = '  -3.29%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.72'
$ws.Cells.Item(33, 5).Value = '  -2.26%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '7.14'
$ws.Cells.Item(34, 5).Value = '  -2.30%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.151'
$ws.Cells.Item(35, 5).Value = '  -8.48%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '38.29'
$ws.Cells.Item(36, 5).Value = '  -8.06%  '
$ws.Cells.Item(37, 5).Value = '  -0.03%  '
$ws.Cells.Item(38, 2).Value = 'OKB'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '54.65'
$ws.Cells.Item(38, 5).Value = '  -4.68%  '
$ws.Cells.Item(39, 2).Value = 'PEPE'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(39, 4).Value = '0.0₃0778'
$ws.Cells.Item(39, 5).Value = '  +7.49%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.08'
$ws.Cells.Item(40, 5).Value = '  -0.88%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0456'
$ws.Cells.Item(41, 5).Value = '  -7.29%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.00'
$ws.Cells.Item(42, 5).Value = '  +0.33%  '
$ws.Cells.Item(43, 2).Value = 'Monero'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '149.74'
$ws.Cells.Item(43, 5).Value = '  +0.74%  '
$ws.Cells.Item(44, 2).Value = 'Stellar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.136'
$ws.Cells.Item(44, 5).Value = '  -9.27%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.48'
$ws.Cells.Item(45, 5).Value = '  +2.09%  '
$ws.Cells.Item(46, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.14'
$ws.Cells.Item(46, 5).Value = '  -4.78%  '
$ws.Cells.Item(47, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.30'
$ws.Cells.Item(47, 5).Value = '  -3.99%  '
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '26.39'
$ws.Cells.Item(48, 5).Value = '  -11.14%  '
$ws.Cells.Item(49, 2).Value = 'ARBITRUM'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.07'
$ws.Cells.Item(49, 5).Value = '  -2.16%  '
$ws.Cells.Item(50, 2).Value = 'WEMIXToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.56'
$ws.Cells.Item(50, 5).Value = '  -3.71%  '
$ws.Cells.Item(51, 2).Value = 'Stacks'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.78'
$ws.Cells.Item(51, 5).Value = '  -4.00%  '
